# edit.ps1 - applies the session_001.docx diff via Word COM-interop.
#
# Strategy: the document is a flat list of top-level paragraphs. We walk the
# edits from the BOTTOM of the document upward so that paragraph indices
# computed up-front stay valid as we mutate the document (only a single new
# paragraph is inserted, and it is inserted last in our bottom-up walk so it
# never perturbs the index of any paragraph we still need to touch).
#
# Helper: replace the visible text of paragraph $index with $newText while
# leaving the paragraph's own run/paragraph formatting (style, rPr, pPr)
# completely untouched - we only rewrite the characters up to (but not
# including) the trailing paragraph mark.
function Set-ParaText($doc, $index, $newText) {
    $para = $doc.Paragraphs.Item($index)
    $full = $para.Range
    $body = $doc.Range($full.Start, $full.End - 1)
    $body.Text = $newText
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Activity Block 4 (old) content, paragraphs 25-28 -> becomes the new
# Activity Block 4 content (Backhand Volley Straight Drop - Drive /
# Boast-Cross-Drive With Counter Drops block).
# ---------------------------------------------------------------------
Set-ParaText $d 25 "3 min: Drill: Drop-Drive: Volley Straight Drop - Drive (Deep Only) (Backhand)"
Set-ParaText $d 26 "(Rule: Drive: The first bounce of all drives must land behind the T-line....)"
Set-ParaText $d 27 "3 min: Drill: Boast-Cross-Drive: Boast-Cross-Drive With Counter Drops (Backhand)"
Set-ParaText $d 28 "(Rule: Cross & Drive: Every cross-court shot and the subsequent drive...)"

# ---------------------------------------------------------------------
# Activity Block 3 (old) content, paragraphs 19-22 -> becomes the new
# Activity Block 3 content (Forehand Volley Straight Drop - Drive /
# Boast-Cross-Drive With Counter Drops block).
# ---------------------------------------------------------------------
Set-ParaText $d 19 "3 min: Drill: Drop-Drive: Volley Straight Drop - Drive (Deep Only) (Forehand)"
Set-ParaText $d 20 "(Rule: Drive: The first bounce of all drives must land behind the T-line....)"
Set-ParaText $d 21 "3 min: Drill: Boast-Cross-Drive: Boast-Cross-Drive With Counter Drops (Forehand)"
Set-ParaText $d 22 "(Rule: Cross & Drive: Every cross-court shot and the subsequent drive...)"

# ---------------------------------------------------------------------
# Activity Block 2 content, paragraphs 14-16: gains a 4th paragraph
# (a "Rule" line) that did not exist before, so insert it first (after
# paragraph 16) and then fix up the text of 14/16 plus the brand new one.
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertParagraphAfter()

$newRule = $d.Paragraphs.Item(17)
$newRule.Style = "Normal"
$newRule.Format.LeftIndent = 36
$newRule.Format.SpaceAfter = 6
$newBody = $d.Range($newRule.Range.Start, $newRule.Range.End - 1)
$newBody.Text = "(Rule: All shots excluding the boast must land behind the T-line.)"
$newBody.Font.Italic = $true
$newBody.Font.Color = 8421504

Set-ParaText $d 16 "11 pts: Conditioned Game: Boast-Cross-Drive: Boast-Cross-Drive Deep Only (Forehand)"
Set-ParaText $d 15 "(Rule: Straight Lob: Must be hit above the service line on the front wall,...)"
Set-ParaText $d 14 "9 pts: Conditioned Game: Drop-Drive: Volley Straight Drop - Straight Lob (Forehand)"

# ---------------------------------------------------------------------
# Activity Block 1 content, paragraphs 8-11.
# ---------------------------------------------------------------------
Set-ParaText $d 8  "9 pts: Conditioned Game: Drop-Drive: Volley Straight Drop - Straight Lob (Backhand)"
Set-ParaText $d 9  "(Rule: Straight Lob: Must be hit above the service line on the front wall,...)"
Set-ParaText $d 10 "11 pts: Conditioned Game: Boast-Cross-Drive: Boast-Cross-Drive Deep Only (Backhand)"
Set-ParaText $d 11 "(Rule: All shots excluding the boast must land behind the T-line.)"

# ---------------------------------------------------------------------
# Warm-up drill text tweak.
# ---------------------------------------------------------------------
Set-ParaText $d 5 "3 min: Drill: Warmup: Compound Boast-Drive-Drop-Drive + 2 shots"

# ---------------------------------------------------------------------
# Session Focus line: only the second run (after the bold "Session
# Focus:" label) changes, so scope the Find/Replace to that paragraph.
# ---------------------------------------------------------------------
$focusPara = $d.Paragraphs.Item(3)
$focusRange = $focusPara.Range
$focusRange.Find.Execute(" Progressive ShoteSide (Archetype: Progressive ShoteSide)", $true, $false, $false, $false, $false, $true, 1, $false, " Dynamic Block Session (Archetype: Dynamic Block Session)", 2)
